$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update F2: convert formula-based hyperlink to a plain text URL
$ws.Range("F2").Value = 'https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/28371/?location=India&locationId=300000000440677&locationLevel=country&mode=location'

# 2) Insert 15 new rows at row 10 (pushes existing rows 10-63 down to 25-78)
$ws.Rows.Item(10).Resize(15).Insert()

# 3) Remove the 9 trailing duplicate rows that are now at 70-78
#    (these were the old rows 55-63, which are superseded by the new rows 10-24)
$ws.Rows.Item(70).Resize(9).Delete()

# 4) Populate the 15 newly inserted rows (10-24) with their data
$ws.Rows.Item(10).RowHeight = 80
$ws.Range("A10").Value = 'DevelopmentAid'
$ws.Range("B10").Value = 'Support in the Implementation of Basin Management Measures and District Ganga Plans'
$ws.Range("D10").Value = 'Governance, Learning'
$ws.Range("F10").Value = 'https://www.developmentaid.org/tenders/view/1550511/support-in-the-implementation-of-basin-management-measures-and-district-ganga-plans'

$ws.Rows.Item(11).RowHeight = 80
$ws.Range("A11").Value = 'DevelopmentAid'
$ws.Range("B11").Value = 'EoI - For Empanelment of agencies for Impact Assessment / Evaluation Services under Programmes for Development and Impact'
$ws.Range("D11").Value = 'Governance, Learning'
$ws.Range("F11").Value = 'https://www.developmentaid.org/tenders/view/1609663/eoi-for-empanelment-of-agencies-for-impact-assessment-evaluation-services-under-programmes-for-devel'

$ws.Rows.Item(12).RowHeight = 80
$ws.Range("A12").Value = 'DevelopmentAid'
$ws.Range("B12").Value = 'RFP - Early-Stage Impact Assessment-CDRI SWP 23-26'
$ws.Range("D12").Value = 'Governance'
$ws.Range("F12").Value = 'https://www.developmentaid.org/tenders/view/1609606/rfp-early-stage-impact-assessment-cdri-swp-23-26'

$ws.Rows.Item(13).RowHeight = 80
$ws.Range("A13").Value = 'DevelopmentAid'
$ws.Range("B13").Value = 'Mirova Energy Transition Emerging Asia (METEA)'
$ws.Range("D13").Value = 'Climate'
$ws.Range("F13").Value = 'https://www.developmentaid.org/tenders/view/1609692/mirova-energy-transition-emerging-asia-metea'

$ws.Rows.Item(14).RowHeight = 80
$ws.Range("A14").Value = 'DevelopmentAid'
$ws.Range("B14").Value = 'Ausbau Erneuerbarer Energien / Expansion of renewable energies'
$ws.Range("D14").Value = 'Climate'
$ws.Range("F14").Value = 'https://www.developmentaid.org/tenders/view/1609597/ausbau-erneuerbarer-energien-expansion-of-renewable-energies'

$ws.Rows.Item(15).RowHeight = 80
$ws.Range("A15").Value = 'DevelopmentAid'
$ws.Range("B15").Value = 'Hiring consultancy firm for filing GST, VAT, IGST Refund Claims and related tasks for the ILO New Delhi office'
$ws.Range("D15").Value = 'Governance'
$ws.Range("F15").Value = 'https://www.developmentaid.org/tenders/view/1609700/hiring-consultancy-firm-for-filing-gst-vat-igst-refund-claims-and-related-tasks-for-the-ilo-new-delh'

$ws.Rows.Item(16).RowHeight = 80
$ws.Range("A16").Value = 'DevelopmentAid'
$ws.Range("B16").Value = 'LOAN-4411 IND: Building India''s Clean Plant Program - IARI-staff-C S09 IARI staff - virologists and other staff - Lab Assistant (02) (Plant Material Management) (57041-001)'
$ws.Range("D16").Value = 'Governance'
$ws.Range("F16").Value = 'https://www.developmentaid.org/tenders/view/1602442/loan-4411-ind-building-indias-clean-plant-program-iari-staff-c-s09-iari-staff-virologists-and-other'

$ws.Rows.Item(17).RowHeight = 80
$ws.Range("A17").Value = 'DevelopmentAid'
$ws.Range("B17").Value = 'LOAN-4291 IND: Himachal Pradesh Subtropical Horticulture, Irrigation, and Value Addition Project - CS 10 Communications and digital contents for Information Education and Knowledge Management under HPSHIVA Project (53189-002)'
$ws.Range("D17").Value = 'Governance, Learning'
$ws.Range("F17").Value = 'https://www.developmentaid.org/tenders/view/1609118/loan-4291-ind-himachal-pradesh-subtropical-horticulture-irrigation-and-value-addition-project-cs-10'

$ws.Rows.Item(18).RowHeight = 80
$ws.Range("A18").Value = 'DevelopmentAid'
$ws.Range("B18").Value = 'TA-10446 REG: Improving Governance for Urban and Water Service Providers - Knowledge and Learning Specialist (National) (58369-001)'
$ws.Range("D18").Value = 'Governance, Learning'
$ws.Range("F18").Value = 'https://www.developmentaid.org/tenders/view/1609114/ta-10446-reg-improving-governance-for-urban-and-water-service-providers-knowledge-and-learning-speci'

$ws.Rows.Item(19).RowHeight = 80
$ws.Range("A19").Value = 'DevelopmentAid'
$ws.Range("B19").Value = 'TA-10446 REG: Improving Governance for Urban and Water Service Providers - Fecal Sludge Management and Urban Institutional Expert (National) (58369-001)'
$ws.Range("D19").Value = 'Governance'
$ws.Range("F19").Value = 'https://www.developmentaid.org/tenders/view/1609113/ta-10446-reg-improving-governance-for-urban-and-water-service-providers-fecal-sludge-management-and'

$ws.Rows.Item(20).RowHeight = 80
$ws.Range("A20").Value = 'DevelopmentAid'
$ws.Range("B20").Value = 'Request For Proposal For Hiring of Technical Partner Agency for Designing, Development and Implementation of AI-Enabled Learning and Knowledge Sharing Platform for ISA'
$ws.Range("D20").Value = 'Governance, Learning'
$ws.Range("F20").Value = 'https://www.developmentaid.org/tenders/view/1589295/request-for-proposal-for-hiring-of-technical-partner-agency-for-designing-development-and-implementa'

$ws.Rows.Item(21).RowHeight = 80
$ws.Range("A21").Value = 'DevelopmentAid'
$ws.Range("B21").Value = 'P168633- Kerala Solid Waste Management Project (ULB PGT PALAKKAD)'
$ws.Range("D21").Value = 'Governance, Climate'
$ws.Range("F21").Value = 'https://www.developmentaid.org/tenders/view/1095526/india-south-asia-p168633-kerala-solid-waste-management-project-ulb-pgt-palakkad-procurement-plan'

$ws.Rows.Item(22).RowHeight = 80
$ws.Range("A22").Value = 'DevelopmentAid'
$ws.Range("B22").Value = 'P166020 - West Bengal Transport and Logistics Spatial Development Project - Procurement Plan (West Bengal Transport Infrastructure Development Corporation Limited)'
$ws.Range("D22").Value = 'Learning, Safety'
$ws.Range("F22").Value = 'https://www.developmentaid.org/tenders/view/485630/india-south-asia-p166020-west-bengal-transport-and-logistics-spatial-development-project-procurement'

$ws.Rows.Item(23).RowHeight = 80
$ws.Range("A23").Value = 'DevelopmentAid'
$ws.Range("B23").Value = 'P179935- Enhancing Landscape and Ecosystem Management (ELEMENT) Project - Procurement Plan (State Forest Development Agency, Tripura)'
$ws.Range("D23").Value = 'Governance, Learning'
$ws.Range("F23").Value = 'https://www.developmentaid.org/tenders/view/1161054/india-south-asia-p179935-enhancing-landscape-and-ecosystem-management-project-procurement-plan-engli'

$ws.Rows.Item(24).RowHeight = 80
$ws.Range("A24").Value = 'DevelopmentAid'
$ws.Range("B24").Value = 'Request for Proposals (RFP): Consulting Agency for Identifying High Growth Livelihoods and Market Access Challenges'
$ws.Range("D24").Value = 'Governance'
$ws.Range("F24").Value = 'https://www.developmentaid.org/tenders/view/1598904/request-for-proposals-rfp-consulting-agency-for-identifying-high-growth-livelihoods-and-market-acces'

# 5) Convert the formula-based Apply_Link formulas in the shifted rows (25-69) to plain text URLs
$ws.Range("F25").Value = 'https://www.developmentaid.org/tenders/view/1602444/loan-4411-ind-building-indias-clean-plant-program-iari-staff-c-s09-iari-staff-virologists-and-other'
$ws.Range("F26").Value = 'https://www.developmentaid.org/tenders/view/1609019/ausbau-erneuerbarer-energien-expansion-of-renewable-energies'
$ws.Range("F27").Value = 'https://www.developmentaid.org/tenders/view/1608947/integrated-urban-climate-action-for-low-carbon-resilient-cities-urban-act'
$ws.Range("F28").Value = 'https://www.developmentaid.org/tenders/view/1608977/integrated-project-for-source-sustainability-and-climate-resilient-rain-fed-agriculture-in-himachal'
$ws.Range("F29").Value = 'https://www.developmentaid.org/tenders/view/1579382/integrated-project-for-source-sustainability-and-climate-resilient-rain-fed-agriculture-in-himachal'
$ws.Range("F30").Value = 'https://www.developmentaid.org/tenders/view/1608441/ta-9849-ind-ta-9849-ind-india-urban-and-water-projects-support-facility-53067-002-environment-safegu'
$ws.Range("F31").Value = 'https://www.developmentaid.org/tenders/view/1608429/ta-10721-ind-support-for-development-of-sustainable-metro-rail-projects-for-nagpur-metro-rail-projec'
$ws.Range("F32").Value = 'https://www.developmentaid.org/tenders/view/1593414/81322849-capacity-development-and-organisational-support-for-convergent-implementation-of-measures-t'
$ws.Range("F33").Value = 'https://www.developmentaid.org/tenders/view/1608627/consulting-feasibility-study-climate-resilience-and-ecosystem-services-in-forest-landscapes'
$ws.Range("F34").Value = 'https://www.developmentaid.org/tenders/view/1605060/loan-4322-ind-enhancing-connectivity-and-sustainability-in-bihar-roads-project-1-expression-of-inter'
$ws.Range("F35").Value = 'https://www.developmentaid.org/tenders/view/1603067/ta-10469-ind-knowledge-and-capacity-building-for-catalyzing-green-growth-and-strengthening-climate-r'
$ws.Range("F36").Value = 'https://www.developmentaid.org/tenders/view/1608411/expression-of-interest-for-selection-of-business-associates-partners-jv-consortium-partnership-for-e'
$ws.Range("F37").Value = 'https://www.developmentaid.org/tenders/view/1608368/43253-025-karnataka-integrated-urban-water-management-investment-program-tranche-1'
$ws.Range("F38").Value = 'https://www.developmentaid.org/tenders/view/1606067/punjab-municipal-services-improvement-project-direct-rfp-for-hiring-of-agency-for-communication-outr'
$ws.Range("F39").Value = 'https://www.developmentaid.org/tenders/view/1411799/46166-003-supporting-human-capital-development-in-meghalaya-phase-2-project-procurement-plan'
$ws.Range("F40").Value = 'https://www.developmentaid.org/tenders/view/1608266/construction-of-swiss-pavilion-for-ai-summit-as-per-the-design-and-space-provided-at-pragati-maidan-'
$ws.Range("F41").Value = 'https://www.developmentaid.org/tenders/view/435581/india-south-asia-p168310-state-of-maharashtras-agribusiness-and-rural-transformation-project-procure'
$ws.Range("F42").Value = 'https://www.developmentaid.org/tenders/view/1607893/ta-6822-ind-support-for-strengthening-multimodal-and-integrated-logistics-ecosystem-development-of-s'
$ws.Range("F43").Value = 'https://www.developmentaid.org/tenders/view/1607892/ta-6822-ind-support-for-strengthening-multimodal-and-integrated-logistics-ecosystem-development-of-s'
$ws.Range("F44").Value = 'https://www.developmentaid.org/tenders/view/1607886/ta-9950-reg-pf-data-science-expert-pf-data-scientist-pf-data-science-expert-54079-001'
$ws.Range("F45").Value = 'https://www.developmentaid.org/tenders/view/1607884/loan-4623-ind-assam-urban-sector-development-project-ausdpcscom-community-mobilization-consultant-cm'
$ws.Range("F46").Value = 'https://www.developmentaid.org/tenders/view/1605066/ta-10488-reg-city-resilience-affordable-housing-sustainable-tourism-and-inclusive-economic-growth-fo'
$ws.Range("F47").Value = 'https://www.developmentaid.org/tenders/view/1608036/hiring-of-event-management-agency-for-organisation-of-two-workshops-for-the-ilo-in-february-2026-in'
$ws.Range("F48").Value = 'https://www.developmentaid.org/tenders/view/1335085/punjab-outcomes-acceleration-in-school-education-operation-p500564'
$ws.Range("F49").Value = 'https://www.developmentaid.org/tenders/view/1607773/strengthening-coastal-resilience-and-the-economy-project-consultancy-service-for-supervision-of-cons'
$ws.Range("F50").Value = 'https://www.developmentaid.org/tenders/view/1607851/rfp-for-for-developing-integrated-one-health-portal-for-enhanced-collaborative-s'
$ws.Range("F51").Value = 'https://www.developmentaid.org/tenders/view/1531583/implementation-and-operation-of-flr-monitoring-evaluation-and-reporting-frameworks-for-forest-landsc'
$ws.Range("F52").Value = 'https://www.developmentaid.org/tenders/view/1534310/81319406-engagement-of-agencyconsortium-for-the-global-best-practices-and-climate-risk-finance-for-h'
$ws.Range("F53").Value = 'https://www.developmentaid.org/tenders/view/1528202/renewable-energy-policy-advisory-services-and-capacity-development-of-government-officials-in-indian'
$ws.Range("F54").Value = 'https://www.developmentaid.org/tenders/view/1607814/unterstutzung-fur-die-partnerschaft-fur-grune-und-nachhaltige-entwicklung-support-for-the-partnershi'
$ws.Range("F55").Value = 'https://www.developmentaid.org/tenders/view/1607571/rfp-river-health-and-water-dependency-assessment-betwa-river-basin'
$ws.Range("F56").Value = 'https://www.developmentaid.org/tenders/view/539175/india-south-asia-p172213-nagaland-enhancing-classroom-teaching-and-resources-procurement-plan'
$ws.Range("F57").Value = 'https://www.developmentaid.org/tenders/view/1607502/ta-6822-ind-support-for-strengthening-multimodal-and-integrated-logistics-ecosystem-development-of-s'
$ws.Range("F58").Value = 'https://www.developmentaid.org/tenders/view/1583551/final-evaluation-of-the-market-development-facility-phase-2-asia-pacific'
$ws.Range("F59").Value = 'https://www.developmentaid.org/tenders/view/1607171/integrated-urban-climate-action-for-low-carbon-resilient-cities-urban-act'
$ws.Range("F60").Value = 'https://www.developmentaid.org/tenders/view/1607022/rfp-audio-quality-annotation-of-student-voice-recordings'
$ws.Range("F61").Value = 'https://www.developmentaid.org/tenders/view/1607202/unicef-chhattisgarh-is-seeking-manufacturers-for-the-development-and-production-of-an-innovative-fee'
$ws.Range("F62").Value = 'https://www.developmentaid.org/tenders/view/978729/india-south-asia-p178418-tripura-rural-economic-growth-and-service-delivery-project-public-works-dep'
$ws.Range("F63").Value = 'https://www.developmentaid.org/tenders/view/1016520/india-south-asia-p178418-tripura-rural-economic-growth-and-service-delivery-project-samagra-shiksha'
$ws.Range("F64").Value = 'https://www.developmentaid.org/tenders/view/1584610/37909-045-trade-and-supply-chain-finance-capacity-development-technical-assistance'
$ws.Range("F65").Value = 'https://www.developmentaid.org/tenders/view/1607091/meghalaya-multisectoral-project-for-adolescent-wellbeing-empowerment-and-resilience-mpower-gd-8procu'
$ws.Range("F66").Value = 'https://www.developmentaid.org/tenders/view/577045/india-south-asia-p154990-jhelum-and-tawi-flood-recovery-project-jammu-and-kashmir-medical-supplies-c'
$ws.Range("F67").Value = 'https://www.developmentaid.org/tenders/view/1606567/ta-6822-ind-support-for-strengthening-multimodal-and-integrated-logistics-ecosystem-development-of-s'
$ws.Range("F68").Value = 'https://www.developmentaid.org/tenders/view/1606912/rfp-supply-installation-testing-commissioning-and-comprehensive-maintenance-of-solar-energy-solution'
$ws.Range("F69").Value = 'https://www.developmentaid.org/tenders/view/1606542/ta-9970-ind-maharashtra-rural-high-voltage-distribution-system-expansion-program-social-sector-exper'
